$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9546840190887451
$ws.Range("B1").Value = 1.864270687103271
$ws.Range("C1").Value = 4.727403163909912
$ws.Range("D1").Value = 2.576142549514771
$ws.Range("E1").Value = 0.46656534075737
